$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Ratio" column (E): header + D/C formula for every data row ---
$ws.Range("E1").Value = "Ratio"
$ws.Range("E1").Font.Bold = $true

$ws.Range("E2").Formula = "=D2/C2"
$ws.Range("E3:E31").Formula = "=D3/C3"

# Apply the 2-decimal number format used by the author's new "Ratio" column
$ws.Range("E2:E31").NumberFormat = "0.00"

# --- "Average gain" summary cells at the end of each 10-row block ---
$ws.Range("G11").Value = "Average gain"
$ws.Range("G11").Font.Bold = $true
$ws.Range("H11").Formula = "=AVERAGE(E2:E11)"
$ws.Range("H11").NumberFormat = "0.00"

$ws.Range("G21").Value = "Average gain"
$ws.Range("G21").Font.Bold = $true
$ws.Range("H21").Formula = "=AVERAGE(E12:E21)"
$ws.Range("H21").NumberFormat = "0.00"

$ws.Range("G31").Value = "Average gain"
$ws.Range("G31").Font.Bold = $true
$ws.Range("H31").Formula = "=AVERAGE(E22:E31)"
$ws.Range("H31").NumberFormat = "0.00"

# --- Column sizing: hide B, widen G for the new "Average gain" label ---
$ws.Columns("B").ColumnWidth = -0.8333333333333334
$ws.Columns("B").Hidden = $true
$ws.Columns("G").ColumnWidth = 30.385416666666668

# --- Final selection, matching the saved workbook state ---
$ws.Range("M31").Select() | Out-Null
